# Apply the "Beginning replacement of mapping from procedure name to
# procedure group" edit to the NotificationPlacingWayName worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to the new machine-readable names ---
$ws.Cells.Item(1, 1).Value = "NotificationPlacingWayName"
$ws.Cells.Item(1, 2).Value = "Freq"
$ws.Cells.Item(1, 3).Value = "TenderProcedureGroup"
$ws.Cells.Item(1, 4).Value = "TenderProcedureDiscretion"

# --- Row 2: discretion label "NA" -> "Other" ---
$ws.Cells.Item(2, 4).Value = "Other"

# --- Rows 11 & 12: swap which procedure name goes with which Freq,
#     keeping group/discretion ("Request for quotes" / "Higher discretion")
#     the same on both rows ---
$a11 = $ws.Cells.Item(11, 1).Value2
$a12 = $ws.Cells.Item(12, 1).Value2
$b11 = $ws.Cells.Item(11, 2).Value2
$b12 = $ws.Cells.Item(12, 2).Value2

$ws.Cells.Item(11, 1).Value = $a12
$ws.Cells.Item(12, 1).Value = $a11
$ws.Cells.Item(11, 2).Value = $b12
$ws.Cells.Item(12, 2).Value = $b11

# Row 11 reverts to default (auto) height, row 12 becomes the tall
# (ht=65) row that row 11 used to be.
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).RowHeight = 65

# --- New rows 15 & 16: two more procedure-name -> group/discretion rows ---
$ws.Cells.Item(15, 1).Value = "Предварительный отбор"
$ws.Cells.Item(15, 2).Value = 1
$ws.Cells.Item(15, 3).Value = "Preliminary selection"
$ws.Cells.Item(15, 4).Value = "Other"

$ws.Cells.Item(16, 1).Value = "Сообщение о заинтересованности в проведении открытого конкурса"
$ws.Cells.Item(16, 2).Value = 1
$ws.Cells.Item(16, 3).Value = "Registration of interest in open tender"
$ws.Cells.Item(16, 4).Value = "Other"

# --- Comment on A1 documenting the old column name ---
$comment = $ws.Range("A1").AddComment()
$comment.Text("Shaun McGirr:" + [char]10 + "Later recoded to ""Type of procedure""")
